# PROS-9738 - CCRU - New POS 2019 KPIs
#
# The "Canteen" channel column is being unified to a single "Canteen_EDU"
# value for every data row (it previously mixed "Canteen" and, on one row,
# the legacy "Canteen_TRAD" label). The view is also scrolled back to the
# top of the sheet, and a fresh copy of the (hidden) AutoFilter database
# defined name is registered, mirroring the workbook's existing
# _xlnm._FilterDatabase_0 / _0_0 / _0_0_0 chain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Normalise the whole "Channel" column (C2:C57 - every data row below the
#    header) to "Canteen_EDU". This covers both the plain "Canteen" cells
#    and the one-off "Canteen_TRAD" cell on row 43.
$ws.Range("C2:C57").Value = "Canteen_EDU"

# 2. Reset the view: scroll/select back to A2 (just below the frozen header
#    row) instead of the previous mid-sheet selection.
[void]$ws.Range("A2").Select()

# 3. Register another generation of the AutoFilter's defined name, matching
#    the existing _xlnm._FilterDatabase_0_0_0 -> _xlnm._FilterDatabase_0_0_0_0
#    progression already present in the workbook.
[void]$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0", "=Canteen!`$A`$1:`$AL`$1")
